$wb = $excel.ActiveWorkbook

# --- Rename sheets ---
$wsOriginal = $wb.Worksheets.Item("original")
$wsOriginal.Name = "Original"

$wsAug = $wb.Worksheets.Item("data augmentation")
$wsAug.Name = "Data augmented"

# --- Sheet "Original": new header row with A1 "Site (4,5,6)" ---
$wsOriginal.Range("A1").Value = "Site (4,5,6)"
$wsOriginal.Range("B1").Value = "Barriers (eV)"
$wsOriginal.Range("C1").Value = "Gap (eV)"
